$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Gdnf"
$ws.Cells.Item(2, 3).Value = "Ret"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02098366666666666
$ws.Cells.Item(2, 8).Value = 0.06295099999999999
$ws.Cells.Item(2, 9).Value = 0.04042238960271747
$ws.Cells.Item(2, 10).Value = 0.04042238960271747
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.328495333333334
$ws.Cells.Item(2, 14).Value = 12.985486
$ws.Cells.Item(2, 15).Value = 0.2325244328639614
$ws.Cells.Item(2, 16).Value = 0.2325244328639614
$ws.Cells.Item(2, 17).Value = 0.09082770324288889
$ws.Cells.Item(2, 18).Value = 0.8174493291859999
$ws.Cells.Item(2, 19).Value = 0.009399193217377968
$ws.Cells.Item(2, 20).Value = 0.009399193217377968

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Gdnf"
$ws.Cells.Item(3, 3).Value = "Ret"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02098366666666666
$ws.Cells.Item(3, 8).Value = 0.06295099999999999
$ws.Cells.Item(3, 9).Value = 0.04042238960271747
$ws.Cells.Item(3, 10).Value = 0.04042238960271747
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.361539333333333
$ws.Cells.Item(3, 14).Value = 22.084618
$ws.Cells.Item(3, 15).Value = 0.3954579193622197
$ws.Cells.Item(3, 16).Value = 0.3954579193622196
$ws.Cells.Item(3, 17).Value = 0.1544720875242222
$ws.Cells.Item(3, 18).Value = 1.390248787718
$ws.Cells.Item(3, 19).Value = 0.01598535408793967
$ws.Cells.Item(3, 20).Value = 0.01598535408793967

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gdnf"
$ws.Cells.Item(4, 3).Value = "Ret"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02098366666666666
$ws.Cells.Item(4, 8).Value = 0.06295099999999999
$ws.Cells.Item(4, 9).Value = 0.04042238960271747
$ws.Cells.Item(4, 10).Value = 0.04042238960271747
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.912382333333333
$ws.Cells.Item(4, 14).Value = 20.737147
$ws.Cells.Item(4, 15).Value = 0.3713294477689628
$ws.Cells.Item(4, 16).Value = 0.3713294477689628
$ws.Cells.Item(4, 17).Value = 0.1450471267552222
$ws.Cells.Item(4, 18).Value = 1.305424140797
$ws.Cells.Item(4, 19).Value = 0.01501002360867894
$ws.Cells.Item(4, 20).Value = 0.01501002360867894

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gdnf"
$ws.Cells.Item(5, 3).Value = "Ret"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02098366666666666
$ws.Cells.Item(5, 8).Value = 0.06295099999999999
$ws.Cells.Item(5, 9).Value = 0.04042238960271747
$ws.Cells.Item(5, 10).Value = 0.04042238960271747
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.012811
$ws.Cells.Item(5, 14).Value = 0.038433
$ws.Cells.Item(5, 15).Value = 0.0006882000048562394
$ws.Cells.Item(5, 16).Value = 0.0006882000048562392
$ws.Cells.Item(5, 17).Value = 0.0002688217536666667
$ws.Cells.Item(5, 18).Value = 0.002419395783
$ws.Cells.Item(5, 19).Value = 0.00002781868872089096
$ws.Cells.Item(5, 20).Value = 0.00002781868872089096

$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Gdnf"
$ws.Cells.Item(6, 3).Value = "Ret"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.4981263333333333
$ws.Cells.Item(6, 8).Value = 1.494379
$ws.Cells.Item(6, 9).Value = 0.9595776103972825
$ws.Cells.Item(6, 10).Value = 0.9595776103972825
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.328495333333334
$ws.Cells.Item(6, 14).Value = 12.985486
$ws.Cells.Item(6, 15).Value = 0.2325244328639614
$ws.Cells.Item(6, 16).Value = 0.2325244328639614
$ws.Cells.Item(6, 17).Value = 2.156137509243778
$ws.Cells.Item(6, 18).Value = 19.405237583194
$ws.Cells.Item(6, 19).Value = 0.2231252396465834
$ws.Cells.Item(6, 20).Value = 0.2231252396465834

$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Gdnf"
$ws.Cells.Item(7, 3).Value = "Ret"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.4981263333333333
$ws.Cells.Item(7, 8).Value = 1.494379
$ws.Cells.Item(7, 9).Value = 0.9595776103972825
$ws.Cells.Item(7, 10).Value = 0.9595776103972825
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.361539333333333
$ws.Cells.Item(7, 14).Value = 22.084618
$ws.Cells.Item(7, 15).Value = 0.3954579193622197
$ws.Cells.Item(7, 16).Value = 0.3954579193622196
$ws.Cells.Item(7, 17).Value = 3.666976595802444
$ws.Cells.Item(7, 18).Value = 33.002789362222
$ws.Cells.Item(7, 19).Value = 0.37947256527428
$ws.Cells.Item(7, 20).Value = 0.37947256527428

$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Gdnf"
$ws.Cells.Item(8, 3).Value = "Ret"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.4981263333333333
$ws.Cells.Item(8, 8).Value = 1.494379
$ws.Cells.Item(8, 9).Value = 0.9595776103972825
$ws.Cells.Item(8, 10).Value = 0.9595776103972825
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 6.912382333333333
$ws.Cells.Item(8, 14).Value = 20.737147
$ws.Cells.Item(8, 15).Value = 0.3713294477689628
$ws.Cells.Item(8, 16).Value = 0.3713294477689628
$ws.Cells.Item(8, 17).Value = 3.443239666301444
$ws.Cells.Item(8, 18).Value = 30.989156996713
$ws.Cells.Item(8, 19).Value = 0.3563194241602839
$ws.Cells.Item(8, 20).Value = 0.3563194241602838

$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Gdnf"
$ws.Cells.Item(9, 3).Value = "Ret"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.4981263333333333
$ws.Cells.Item(9, 8).Value = 1.494379
$ws.Cells.Item(9, 9).Value = 0.9595776103972825
$ws.Cells.Item(9, 10).Value = 0.9595776103972825
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.012811
$ws.Cells.Item(9, 14).Value = 0.038433
$ws.Cells.Item(9, 15).Value = 0.0006882000048562394
$ws.Cells.Item(9, 16).Value = 0.0006882000048562392
$ws.Cells.Item(9, 17).Value = 0.006381496456333333
$ws.Cells.Item(9, 18).Value = 0.057433468107
$ws.Cells.Item(9, 19).Value = 0.0006603813161353484
$ws.Cells.Item(9, 20).Value = 0.0006603813161353483
